# Hortaliza, Femacal de La Calera - Brocoli: weekly refresh.
#
# The data block of rows 328:363 (18 date-pairs, "Primera"/"Segunda" rows)
# shifts down by two rows to 330:365, and two brand-new rows (a fresh
# date-pair of observations) are inserted at the top of the block at
# 328:329. The dimension / used range grows from R363 to R365.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at 328, pushing the former 328:363 block
# down to 330:365 (two single-row inserts so both rows actually shift,
# rather than a single two-row range insert).
$ws.Rows.Item(328).Insert()
$ws.Rows.Item(328).Insert()

# Row 328 ("Primera" quality record for the new date).
$ws.Cells.Item(328, 1).Value = 3
$ws.Cells.Item(328, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(328, 3).Value = "Coquimbo"
$ws.Cells.Item(328, 4).Value = 44449
$ws.Cells.Item(328, 5).Value = 5
$ws.Cells.Item(328, 6).Value = 100112023
$ws.Cells.Item(328, 7).Value = "Brócoli"
$ws.Cells.Item(328, 8).Value = "Sin especificar"
$ws.Cells.Item(328, 9).Value = "Primera"
$ws.Cells.Item(328, 10).Value = 3000
$ws.Cells.Item(328, 11).Value = 550
$ws.Cells.Item(328, 12).Value = 600
$ws.Cells.Item(328, 13).Value = 573
$ws.Cells.Item(328, 14).Value = "`$/unidad"
$ws.Cells.Item(328, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(328, 16).Value = 573
$ws.Cells.Item(328, 17).Value = 1
$ws.Cells.Item(328, 18).Value = "Hortaliza"

# Row 329 ("Segunda" quality record for the new date).
$ws.Cells.Item(329, 1).Value = 3
$ws.Cells.Item(329, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(329, 3).Value = "Coquimbo"
$ws.Cells.Item(329, 4).Value = 44449
$ws.Cells.Item(329, 5).Value = 5
$ws.Cells.Item(329, 6).Value = 100112023
$ws.Cells.Item(329, 7).Value = "Brócoli"
$ws.Cells.Item(329, 8).Value = "Sin especificar"
$ws.Cells.Item(329, 9).Value = "Segunda"
$ws.Cells.Item(329, 10).Value = 1500
$ws.Cells.Item(329, 11).Value = 450
$ws.Cells.Item(329, 12).Value = 450
$ws.Cells.Item(329, 13).Value = 450
$ws.Cells.Item(329, 14).Value = "`$/unidad"
$ws.Cells.Item(329, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(329, 16).Value = 450
$ws.Cells.Item(329, 17).Value = 1
$ws.Cells.Item(329, 18).Value = "Hortaliza"
